$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(17.84, 0.36, 20.19210028648376, 6.46875),
    @(17.84, 0.36, 19.50743865966797, 6.84375),
    @(17.84, 0.36, 18.69954299926758, 6.921875),
    @(17.84, 0.36, 13.4420440196991, 11),
    @(17.84, 0.36, 13.99964332580566, 10.203125),
    @(37.37, 0.36, 18.73186135292053, 15.734375),
    @(39.78, 0.16, 9.278669595718384, 6.703125),
    @(39.78, 0.16, 8.955015659332275, 7.125),
    @(39.78, 0.16, 8.963810920715332, 6.8125)
)

$startRow = 36
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
